$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("B7").Value = 2

$ws.Range("B10").Select()
